$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two brand-new weekly price entries are inserted right after the existing
# row 54 (pushing the rest of the "Vega Modelo de Temuco - Rabanito" table
# down by two rows, from A1:R76 to A1:R78).
$ws.Rows("55:56").Insert()

# New row 55
$ws.Cells.Item(55, 1).Value = 10
$ws.Cells.Item(55, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(55, 3).Value = "La Araucanía"
$ws.Cells.Item(55, 4).Value = 44784
$ws.Cells.Item(55, 5).Value = 9
$ws.Cells.Item(55, 6).Value = 300000001
$ws.Cells.Item(55, 7).Value = "Rabanito"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 100
$ws.Cells.Item(55, 11).Value = 7000
$ws.Cells.Item(55, 12).Value = 8000
$ws.Cells.Item(55, 13).Value = 7500
$ws.Cells.Item(55, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(55, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(55, 16).Value = 625
$ws.Cells.Item(55, 17).Value = 12
$ws.Cells.Item(55, 18).Value = "Hortaliza"

# New row 56
$ws.Cells.Item(56, 1).Value = 10
$ws.Cells.Item(56, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(56, 3).Value = "La Araucanía"
$ws.Cells.Item(56, 4).Value = 44784
$ws.Cells.Item(56, 5).Value = 9
$ws.Cells.Item(56, 6).Value = 300000001
$ws.Cells.Item(56, 7).Value = "Rabanito"
$ws.Cells.Item(56, 8).Value = "Sin especificar"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 40
$ws.Cells.Item(56, 11).Value = 5000
$ws.Cells.Item(56, 12).Value = 5000
$ws.Cells.Item(56, 13).Value = 5000
$ws.Cells.Item(56, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(56, 15).Value = "Región Metropolitana"
$ws.Cells.Item(56, 16).Value = 417
$ws.Cells.Item(56, 17).Value = 12
$ws.Cells.Item(56, 18).Value = "Hortaliza"
